$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author added a new TODO item ("Check SunEye Obstructions and both
# Solar Pathfinder file imports", owner Janine, status "Not done") above
# the "Simulation status and warnings dialog box" row. Insert a whole row
# at row 42 so every row/formula below (SUM ranges, later rows 43-74,
# etc.) shifts down automatically, exactly like Excel's own
# Rows(...).Insert would do.
$ws.Rows("42:42").Insert() | Out-Null

# Populate the newly inserted row.
$ws.Range("A42").Value = "Not done"
$ws.Range("B42").Value = "Check SunEye Obstructions and both Solar Pathfinder file imports"
$ws.Range("C42").Value = "Janine"

# Match the saved view state: selection moved to A43 (and Excel drops the
# stale topLeftCell scroll position once the target cell is selected).
$ws.Activate() | Out-Null
$ws.Range("A43").Select() | Out-Null
